# Initial commit v2p1, R2020b
#
# Updates the brake-pedal-abstraction workbook:
#  - H13 / H22 on the three "PedalAbstract" data sheets now compute the
#    response-time factor from 1/0.025 (was 1/0.2) -> 40 (was 5).
#  - The active tab moves from "None" to "Bus_Makhulu", and each sheet's
#    last-used cell selection is updated to match.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sedan_HambaLG", "Sedan_Hamba", "Bus_Makhulu")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H13").Formula = "=1/0.025"
    $ws.Range("H22").Formula = "=1/0.025"
}

# Re-point the last selected cell on each sheet (mirrors the saved
# sheetView/selection state in the workbook), then leave "Bus_Makhulu"
# as the active tab/sheet, matching the new activeTab in bookViews.
$ws1 = $wb.Worksheets.Item("Sedan_HambaLG")
$ws1.Activate()
$ws1.Range("H22").Select()

$ws2 = $wb.Worksheets.Item("Sedan_Hamba")
$ws2.Activate()
$ws2.Range("H22").Select()

$ws3 = $wb.Worksheets.Item("Bus_Makhulu")
$ws3.Activate()
$ws3.Range("J22").Select()
